$wb = $excel.ActiveWorkbook

# --- Rename sheets (task-order tab names) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509960951654131"
$wb.Worksheets.Item(2).Name = "NB_TO-16509960986466053"
$wb.Worksheets.Item(3).Name = "RS_TO-16509960986466053"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509960986945407"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509960987585418"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650996095133378.csv"
$ws1.Range("B3").Value = "GNG_stims-16509960951493788.csv"
$ws1.Range("B4").Value = "go_stims-16509960951493788.csv"
$ws1.Range("B5").Value = "GNG_stims-16509960951654131.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_0-16509960958133812.csv"
$ws2.Range("B3").Value = "ZB-match_2-1650996095637415.csv"
$ws2.Range("B4").Value = "OB-1650996097206381.csv"
$ws2.Range("B5").Value = "OB-16509960966694198.csv"
$ws2.Range("B6").Value = "TB-16509960986225467.csv"
$ws2.Range("B7").Value = "OB-16509960974703472.csv"
$ws2.Range("B8").Value = "ZB-match_6-16509960952454154.csv"
$ws2.Range("B9").Value = "TB-16509960981183443.csv"
$ws2.Range("B10").Value = "TB-16509960984225395.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650996098662542.csv"
$ws4.Range("B3").Value = "ZM_stims-16509960986466053.csv"
$ws4.Range("B4").Value = "MM_stims-16509960986785762.csv"
$ws4.Range("B5").Value = "ZM_stims-1650996098662542.csv"
$ws4.Range("B6").Value = "MM_stims-16509960986945407.csv"
$ws4.Range("B7").Value = "ZM_stims-16509960986785762.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509960987425494.csv"
$ws5.Range("B3").Value = "SAT_stims-16509960986945407.csv"
$ws5.Range("B4").Value = "SAT_stims-16509960987105417.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509960987265406.csv"
